# إضافة حدث جديد في Card14
# Row 13 on the Card14 sheet previously had several blank (but present)
# string cells; this event's edit fills them with the literal text "nan"
# (matching the rest of the sheet's placeholder convention), and appends a
# brand-new maintenance-log row (14) underneath it.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Card14")

# --- Row 13: the previously-empty cells now read "nan" ------------------
$ws.Range("B13:K13").Value = "nan"
$ws.Range("M13").Value = "nan"

# --- Row 14: new event row -----------------------------------------------
# Column A holds a card/tone-group id stored as text elsewhere in the
# sheet ("14"), so force text formatting before assigning it, otherwise
# Excel would store it as the number 14.
$ws.Range("A14").NumberFormat = "@"
$ws.Range("A14").Value = "14"

# B14:K14 and M14 stay blank for this event (no tone-range / checkmark
# data recorded), but the cells still need to exist in the row. Touching a
# format property that's already at its default value materialises the
# cell without altering its appearance or introducing a new style.
foreach ($col in @("B","C","D","E","F","G","H","I","J","K","M")) {
    $ws.Range($col + "14").Font.Bold = $false
}

$ws.Range("L14").Value = "13\8\2024"
$ws.Range("N14").Value = "تم تشحيم المكنه بالكامل +عمل صيانه"
$ws.Range("O14").Value = "تيم العمل"
